# Trade #22 closed at 2026-02-16 21:25:41 - momentum DOWN +0.000%
# Append a new trade row (row 5) to the "momentum" sheet, following the
# same shape as the existing rows (2-4).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("momentum")

# Seed row 5 from row 4 (same column layout/types, incl. the two blank
# Exit Price / Exit Reason cells) and then patch only the fields that
# differ for this trade. Doing it this way - rather than writing every
# cell's .Value directly - avoids Excel's autodetection turning the
# "2026-02-16" date-like text into a date serial, and keeps the blank
# G/M cells present (as they are for every other row) instead of being
# dropped the way an explicit empty-string .Value write would be.
$ws.Range("A4:N4").Copy($ws.Range("A5:N5"))

$ws.Range("A5").Value = 22
$ws.Range("C5").Value = "21:25:41"
$ws.Range("F5").Value = 69101.72500000001
$ws.Range("L5").Value = "Downward momentum: -0.262% over 10 samples"
